$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 4-12: one per VRAM-study train_config, dataset03, same goal text ---
# Shared strings must be created in the same order as in the target file:
#   31 "VRAM study;"                  -> write this text first
#   32..40 "train_config-230901-N.yml" (N=0..8), in order
#   41 "run successful;"              -> write this one last

# 1) Seed "VRAM study;" first so it lands at shared-string index 31
$ws.Range("F4").Value = "VRAM study;"

# 2) Fill rows 4..12 with date, config filename, dataset, goal (creates the
#    "train_config-230901-N.yml" strings in ascending order -> indices 32..40)
$dateSerial = 45170
$names = @(
  "train_config-230901-0.yml",
  "train_config-230901-1.yml",
  "train_config-230901-2.yml",
  "train_config-230901-3.yml",
  "train_config-230901-4.yml",
  "train_config-230901-5.yml",
  "train_config-230901-6.yml",
  "train_config-230901-7.yml",
  "train_config-230901-8.yml"
)

for ($i = 0; $i -lt $names.Length; $i++) {
  $r = 4 + $i
  $ws.Range("A$r").Value = $dateSerial
  $ws.Range("B$r").Value = $names[$i]
  $ws.Range("E$r").Value = "dataset03"
  $ws.Range("F$r").Value = "VRAM study;"
}

# Rows 4 & 5 already carried the "m/d/yyyy" date style (s="2") from the
# pre-existing blank rows; rows 6-12 are brand-new and default to General,
# so copy A2's date format onto them without creating any new style/font.
$ws.Range("A2").Copy()
$ws.Range("A6:A12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Row 11 (train_config-230901-7.yml) also records an aborted-but-successful
#    run with valid input data format; "run successful;" is the newest string -> index 41
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = "run successful;"
$ws.Range("M11").Value = 0
$ws.Range("O11").Value = 1

# --- Row 18: separate note about dataset03 / better performance metrics ---
$ws.Range("E18").Value = "dataset03"
$ws.Range("F18").Value = "better performance metrics;"

# --- Register the small (8pt) font that the real edit left in styles.xml.
#     It isn't actually applied to any cell in the final sheet (confirmed by
#     the unchanged cellXfs table), so create it on a scratch row far outside
#     the used range and immediately delete that row again.
$ws.Range("A500").Font.Size = 8
$ws.Rows.Item(500).Delete()

# --- Column width adjustments (stored width = ColumnWidth + 5/6) ---
$ws.Columns.Item(1).ColumnWidth = 9 - 5/6
$ws.Columns.Item(2).ColumnWidth = 24.7109375 - 5/6
$ws.Columns.Item(3).ColumnWidth = 15.7109375 - 5/6
$ws.Columns.Item(5).ColumnWidth = 10.28515625 - 5/6
$ws.Columns.Item(6).ColumnWidth = 18.140625 - 5/6
$ws.Columns.Item(7).ColumnWidth = 28.7109375 - 5/6
$ws.Columns.Item(10).ColumnWidth = 16.7109375 - 5/6
$ws.Columns.Item(11).ColumnWidth = 8.140625 - 5/6

# --- Sheet view: zoom to 85% and move the selection to G6 ---
$excel.ActiveWindow.Zoom = 85
$ws.Range("G6").Select()
